$wb = $excel.ActiveWorkbook

# --- Sheet: 展览 (Exhibition) ---
$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value  = 1665
$ws.Range("F3").Value  = 869
$ws.Range("F4").Value  = 276
$ws.Range("F5").Value  = 85
$ws.Range("F6").Value  = 1195
$ws.Range("F7").Value  = 817
$ws.Range("F8").Value  = 842
$ws.Range("F9").Value  = 1549
$ws.Range("F10").Value = 313
$ws.Range("F11").Value = 1065
$ws.Range("F13").Value = 78
$ws.Range("F16").Value = 521
$ws.Range("F17").Value = 73
$ws.Range("F18").Value = 46
$ws.Range("F19").Value = 13
$ws.Range("F20").Value = 114
$ws.Range("F22").Value = 584
$ws.Range("F23").Value = 589
$ws.Range("F24").Value = 62
$ws.Range("F25").Value = 12
$ws.Range("F26").Value = 786
$ws.Range("F27").Value = 264
$ws.Range("F28").Value = 201

# --- Sheet: 演出 (Performance) ---
$ws = $wb.Worksheets.Item("演出")
$ws.Range("G2").Value  = "不可售"
$ws.Range("F3").Value  = 1043
$ws.Range("F5").Value  = 285
$ws.Range("F7").Value  = 154
$ws.Range("F9").Value  = 600
$ws.Range("F10").Value = 93

# --- Sheet: 本地生活 (Local Life) ---
$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F2").Value = 269

# --- Sheet: 全部类型 (All Types) ---
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value  = 269
$ws.Range("F3").Value  = 1665
$ws.Range("G4").Value  = "不可售"
$ws.Range("F5").Value  = 869
$ws.Range("F6").Value  = 276
$ws.Range("F7").Value  = 1043
$ws.Range("F8").Value  = 85
$ws.Range("F9").Value  = 1195
$ws.Range("F10").Value = 817
$ws.Range("F11").Value = 842
$ws.Range("F12").Value = 1549
$ws.Range("F13").Value = 313
$ws.Range("F14").Value = 1065
$ws.Range("F16").Value = 78
$ws.Range("F19").Value = 521
$ws.Range("F20").Value = 73
$ws.Range("F21").Value = 46
$ws.Range("F23").Value = 13
$ws.Range("F24").Value = 285
$ws.Range("F25").Value = 114
$ws.Range("F28").Value = 154
$ws.Range("F29").Value = 154
$ws.Range("F30").Value = 584
$ws.Range("F31").Value = 589
$ws.Range("F32").Value = 62
$ws.Range("F33").Value = 12
$ws.Range("F34").Value = 786
$ws.Range("F35").Value = 264
$ws.Range("F37").Value = 201
$ws.Range("F38").Value = 600
$ws.Range("F39").Value = 93
$ws.Range("F40").Value = 93
